$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 111114550
$ws.Range("I64").Value = 1000000000
$ws.Range("J64").Value = 3867.5
$ws.Range("K64").Value = 1000000000
$ws.Range("L64").Value = 3867.5
$ws.Range("M64").Value = -999999752
$ws.Range("N64").Value = -4363.5

$ws.Range("H67").Value = 111114550
$ws.Range("I67").Value = 1000000000
$ws.Range("J67").Value = 3867.5
$ws.Range("K67").Value = 1000000000
$ws.Range("L67").Value = 3867.5
$ws.Range("M67").Value = -999999142
$ws.Range("N67").Value = -5583.5

$ws.Range("H74").Value = 3499.25
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3499.25
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3499.25
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -5371.25

$ws.Range("H76").Value = 4645.5884
$ws.Range("I76").Value = 2900
$ws.Range("J76").Value = 6609.375
$ws.Range("K76").Value = 2900
$ws.Range("L76").Value = 6609.375
$ws.Range("M76").Value = -2585
$ws.Range("N76").Value = -7239.375

$ws.Range("H77").Value = 3499.25
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3499.25
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 17496.25
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -26856.25

$ws.Range("H79").Value = 4645.5884
$ws.Range("I79").Value = 2900
$ws.Range("J79").Value = 6609.375
$ws.Range("K79").Value = 2900
$ws.Range("L79").Value = 6609.375
$ws.Range("M79").Value = -1808
$ws.Range("N79").Value = -8793.375

$ws.Range("H98").Value = 2315.3635
$ws.Range("I98").Value = 1949.4286
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 1949.4286
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = -451.4286
$ws.Range("N98").Value = -12996

$ws.Range("H122").Value = 2315.3635
$ws.Range("I122").Value = 1949.4286
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 5848.2858
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -3398.2858
$ws.Range("N122").Value = -34900

$ws.Range("H138").Value = 2668.4255
$ws.Range("I138").Value = 1993.7693
$ws.Range("J138").Value = 3503.7144
$ws.Range("K138").Value = 5981.3079
$ws.Range("L138").Value = 10511.1432
$ws.Range("M138").Value = -841.3078999999998
$ws.Range("N138").Value = -20791.1432

$ws.Range("H141").Value = 5164.9067
$ws.Range("I141").Value = 2874.1282
$ws.Range("J141").Value = 27500
$ws.Range("K141").Value = 8622.384600000001
$ws.Range("L141").Value = 82500
$ws.Range("M141").Value = -3442.384600000001
$ws.Range("N141").Value = -92860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2402.5
$ws.Range("I63").Value = 2405
$ws.Range("J63").Value = 2400
$ws.Range("K63").Value = 2405
$ws.Range("L63").Value = 2400
$ws.Range("M63").Value = -1719
$ws.Range("N63").Value = -3772

$ws.Range("H66").Value = 2402.5
$ws.Range("I66").Value = 2405
$ws.Range("J66").Value = 2400
$ws.Range("K66").Value = 12025
$ws.Range("L66").Value = 12000
$ws.Range("M66").Value = -8593
$ws.Range("N66").Value = -18864

$ws.Range("H88").Value = 3024.3333
$ws.Range("I88").Value = 2261.1667
$ws.Range("J88").Value = 3533.111
$ws.Range("K88").Value = 2261.1667
$ws.Range("L88").Value = 3533.111
$ws.Range("M88").Value = -1855.1667
$ws.Range("N88").Value = -4345.111

$ws.Range("H91").Value = 3024.3333
$ws.Range("I91").Value = 2261.1667
$ws.Range("J91").Value = 3533.111
$ws.Range("K91").Value = 2261.1667
$ws.Range("L91").Value = 3533.111
$ws.Range("M91").Value = -857.1667000000002
$ws.Range("N91").Value = -6341.111

$ws.Range("H132").Value = 5255.49
$ws.Range("I132").Value = 3405.0476
$ws.Range("J132").Value = 13890.889
$ws.Range("K132").Value = 10215.1428
$ws.Range("L132").Value = 41672.667
$ws.Range("M132").Value = -7685.1428
$ws.Range("N132").Value = -46732.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2835.2144
$ws.Range("I86").Value = 2344
$ws.Range("J86").Value = 2998.9524
$ws.Range("K86").Value = 2344
$ws.Range("L86").Value = 2998.9524
$ws.Range("M86").Value = -1221
$ws.Range("N86").Value = -5244.9524

$ws.Range("H89").Value = 2835.2144
$ws.Range("I89").Value = 2344
$ws.Range("J89").Value = 2998.9524
$ws.Range("K89").Value = 11720
$ws.Range("L89").Value = 14994.762
$ws.Range("M89").Value = -6104
$ws.Range("N89").Value = -26226.762

$ws.Range("H105").Value = 1907
$ws.Range("I105").Value = 1373.6364
$ws.Range("K105").Value = 1373.6364
$ws.Range("M105").Value = 373.3635999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4367.7095
$ws.Range("I62").Value = 4895.8335
$ws.Range("J62").Value = 2557
$ws.Range("K62").Value = 4895.8335
$ws.Range("L62").Value = 2557
$ws.Range("M62").Value = -4271.8335
$ws.Range("N62").Value = -3805

$ws.Range("H65").Value = 4367.7095
$ws.Range("I65").Value = 4895.8335
$ws.Range("J65").Value = 2557
$ws.Range("K65").Value = 24479.1675
$ws.Range("L65").Value = 12785
$ws.Range("M65").Value = -21359.1675
$ws.Range("N65").Value = -19025

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1368.9286
$ws.Range("I14").Value = 1368.9286
$ws.Range("K14").Value = 4106.7858
$ws.Range("M14").Value = -3933.7858

$ws.Range("H129").Value = 1191.6428
$ws.Range("I129").Value = 881.6667
$ws.Range("J129").Value = 1424.125
$ws.Range("K129").Value = 2645.0001
$ws.Range("L129").Value = 4272.375
$ws.Range("M129").Value = 2354.9999
$ws.Range("N129").Value = -14272.375

$ws.Range("H131").Value = 905.94116
$ws.Range("I131").Value = 343.75
$ws.Range("J131").Value = 1010.5349
$ws.Range("K131").Value = 1031.25
$ws.Range("L131").Value = 3031.6047
$ws.Range("M131").Value = 4008.75
$ws.Range("N131").Value = -13111.6047

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3914.2856
$ws.Range("I70").Value = 3900
$ws.Range("K70").Value = 3900
$ws.Range("M70").Value = -3630

$ws.Range("H73").Value = 3914.2856
$ws.Range("I73").Value = 3900
$ws.Range("K73").Value = 3900
$ws.Range("M73").Value = -2964

$ws.Range("H80").Value = 9036.666999999999
$ws.Range("I80").Value = 10832.857
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 10832.857
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -9834.857
$ws.Range("N80").Value = -4746

$ws.Range("H83").Value = 9036.666999999999
$ws.Range("I83").Value = 10832.857
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 54164.285
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -49172.285
$ws.Range("N83").Value = -23734

$ws.Range("H101").Value = 42828.5
$ws.Range("J101").Value = 42828.5
$ws.Range("L101").Value = 42828.5
$ws.Range("N101").Value = -49318.5

$ws.Range("H132").Value = 20146.902
$ws.Range("I132").Value = 30205.77
$ws.Range("J132").Value = 2315.2727
$ws.Range("K132").Value = 90617.31
$ws.Range("L132").Value = 6945.8181
$ws.Range("M132").Value = -88087.31
$ws.Range("N132").Value = -12005.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1430.2778
$ws.Range("I82").Value = 1225.8
$ws.Range("J82").Value = 1685.875
$ws.Range("K82").Value = 1225.8
$ws.Range("L82").Value = 1685.875
$ws.Range("M82").Value = -864.8
$ws.Range("N82").Value = -2407.875

$ws.Range("H85").Value = 1430.2778
$ws.Range("I85").Value = 1225.8
$ws.Range("J85").Value = 1685.875
$ws.Range("K85").Value = 1225.8
$ws.Range("L85").Value = 1685.875
$ws.Range("M85").Value = 22.20000000000005
$ws.Range("N85").Value = -4181.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1052.8334
$ws.Range("I122").Value = 986.0909
$ws.Range("J122").Value = 1157.7142
$ws.Range("K122").Value = 2958.2727
$ws.Range("L122").Value = 3473.1426
$ws.Range("M122").Value = -508.2727
$ws.Range("N122").Value = -8373.142599999999

$ws.Range("H123").Value = 21142
$ws.Range("J123").Value = 21429
$ws.Range("L123").Value = 21429
$ws.Range("N123").Value = -31229
